$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add new header cells P1 and Q1 (continuing the 0..15 sequence) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Match the formatting used by the rest of row 1 (bold/border/centered style)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Rows 2-25: update columns I, K, M, O and add new columns P, Q ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
